$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 19460.04417135087
$ws.Range("D3").Value = 905.1407206316592

$ws.Range("B4").Value = 8173.496963332988
$ws.Range("D4").Value = 442.3351594171519

$ws.Range("B5").Value = 6452.864243835624
$ws.Range("D5").Value = 116.3737041095889

$ws.Range("B6").Value = 13341.54491438356
$ws.Range("D6").Value = 457.712813013698

$ws.Range("B7").Value = 16797.49992739726
$ws.Range("D7").Value = 918.6379260273966

$ws.Range("B8").Value = 26967.45173630147
$ws.Range("D8").Value = 1355.441796575342

$ws.Range("B9").Value = 36160.18990958908
$ws.Range("D9").Value = 1417.97261369863

$ws.Range("F10").Value = 28022463.60674525

$ws.Range("G11").Value = 0.8167863866874109

$ws.Range("F12").Value = 1303402.637709589
$ws.Range("G12").Value = 0.04651277831959957

$ws.Range("G13").Value = 0.1367008349929896
